# MISSOURI_2021.xlsx cleanup
# 1. Rename header row (Spanish labels -> clean column codes)
# 2. Apply Proper-Case (Excel PROPER) to the state/municipality name
#    columns (A and B) for all data rows
# 3. Remove the trailing metadata/footer rows (742-746)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -----------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Proper-case columns A and B for data rows 2-740 ----------------
# A helper cell (far outside the used range) is used to evaluate the
# PROPER() worksheet formula for each value and read the result back,
# avoiding self-referencing formulas in A/B themselves.
$helper = $ws.Range("Z1")

for ($r = 2; $r -le 740; $r++) {
    foreach ($col in @("A", "B")) {
        $cell = $ws.Range($col + $r)
        $val = $cell.Value()
        if ($val -ne $null) {
            $helper.Formula = "=PROPER(" + $col + $r + ")"
            $cell.Value = $helper.Value()
        }
    }
}
$helper.Clear()

# --- 3. Drop the trailing metadata rows (742-746) -----------------------
$ws.Range("A742:A746").EntireRow.Delete()
